$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.1204449352092638
$ws.Range("A2").Value = -0.0099999992861405929
$ws.Range("A3").Value = -0.036231432815501563
$ws.Range("A4").Value = 0.28399068916020909
$ws.Range("A5").Value = -0.084096086473598675
$ws.Range("A6").Value = -0.0059999992830803706
$ws.Range("A7").Value = -0.019999999160580373
$ws.Range("A8").Value = -0.019999999153377246
$ws.Range("A9").Value = -0.0059999992653239076
$ws.Range("A10").Value = -0.0059999992585773043
$ws.Range("A11").Value = -0.0044999992714203074
$ws.Range("A12").Value = -0.0059999992571504457
$ws.Range("A13").Value = -0.0059999992517134615
$ws.Range("A14").Value = -0.011999999199057143
$ws.Range("A15").Value = 0.017733910635444161
$ws.Range("A16").Value = -0.0059999992475718855
$ws.Range("A17").Value = -0.0059999992442856254
$ws.Range("A18").Value = -0.0089999992178624311
$ws.Range("A19").Value = -0.0089999992940379414
$ws.Range("A20").Value = -0.0089999992879157276
$ws.Range("A21").Value = -0.0089999992870062329
$ws.Range("A22").Value = -0.0089999992863658562
$ws.Range("A23").Value = -0.0089999992821354624
$ws.Range("A24").Value = -0.041999998987698461
$ws.Range("A25").Value = -0.041999998982228171
$ws.Range("A26").Value = -0.0059999992824835147
$ws.Range("A27").Value = -0.0059999992812675984
$ws.Range("A28").Value = -0.0059999992748469566
$ws.Range("A29").Value = -0.011999999219755253
$ws.Range("A30").Value = -0.019999999150082992
$ws.Range("A31").Value = -0.014999999193213043
$ws.Range("A32").Value = 0.0237147251299481
$ws.Range("A33").Value = -0.0059999992708945626
